# Insert a new price record as row 85 in the daily/weekly logic sheet for
# "Vega Monumental Concepción - Choclo". This pushes the existing rows
# 85-101 down to 86-102 and the sheet dimension grows from R101 to R102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(85).Insert()

$ws.Cells.Item(85, 1).Value  = 11
$ws.Cells.Item(85, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(85, 3).Value  = "Bíobío"
$ws.Cells.Item(85, 4).Value  = 44736
$ws.Cells.Item(85, 5).Value  = 8
$ws.Cells.Item(85, 6).Value  = 100112024
$ws.Cells.Item(85, 7).Value  = "Choclo"
$ws.Cells.Item(85, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(85, 9).Value  = "Primera"
$ws.Cells.Item(85, 10).Value = 100
$ws.Cells.Item(85, 11).Value = 38000
$ws.Cells.Item(85, 12).Value = 40000
$ws.Cells.Item(85, 13).Value = 39000
$ws.Cells.Item(85, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(85, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(85, 16).Value = 557
$ws.Cells.Item(85, 17).Value = 70
$ws.Cells.Item(85, 18).Value = "Hortaliza"
